# Update the "想去人数" (F column) figures across the four worksheets
# (展览 / 演出 / 本地生活 / 全部类型) to match the refreshed crawl output
# committed as "Update gh-pages to output generated at 456a3b4".
#
# Worksheets are addressed by their 1-based index, which corresponds to
# sheet order in the workbook: 1=展览, 2=演出, 3=本地生活, 4=全部类型.

$wb = $excel.ActiveWorkbook

# --- Sheet 1: 展览 ---
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(5, 6).Value = 5091
$ws.Cells.Item(6, 6).Value = 5091
$ws.Cells.Item(7, 6).Value = 99
$ws.Cells.Item(9, 6).Value = 508
$ws.Cells.Item(10, 6).Value = 2
$ws.Cells.Item(11, 6).Value = 1148
$ws.Cells.Item(12, 6).Value = 703
$ws.Cells.Item(13, 6).Value = 4928
$ws.Cells.Item(15, 6).Value = 55
$ws.Cells.Item(16, 6).Value = 73
$ws.Cells.Item(17, 6).Value = 209
$ws.Cells.Item(18, 6).Value = 219
$ws.Cells.Item(19, 6).Value = 97
$ws.Cells.Item(21, 6).Value = 3756
$ws.Cells.Item(24, 6).Value = 3630
$ws.Cells.Item(25, 6).Value = 170
$ws.Cells.Item(26, 6).Value = 166
$ws.Cells.Item(27, 6).Value = 12
$ws.Cells.Item(28, 6).Value = 206
$ws.Cells.Item(29, 6).Value = 231
$ws.Cells.Item(30, 6).Value = 199
$ws.Cells.Item(32, 6).Value = 105
$ws.Cells.Item(36, 6).Value = 6420
$ws.Cells.Item(37, 6).Value = 1016
$ws.Cells.Item(38, 6).Value = 481
$ws.Cells.Item(40, 6).Value = 970
$ws.Cells.Item(42, 6).Value = 1310
$ws.Cells.Item(43, 6).Value = 153
$ws.Cells.Item(44, 6).Value = 638
$ws.Cells.Item(46, 6).Value = 2211
$ws.Cells.Item(47, 6).Value = 310
$ws.Cells.Item(49, 6).Value = 761
$ws.Cells.Item(50, 6).Value = 904

# --- Sheet 2: 演出 ---
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(9, 6).Value = 77
$ws.Cells.Item(20, 6).Value = 48
$ws.Cells.Item(23, 6).Value = 799

# --- Sheet 3: 本地生活 ---
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(2, 6).Value = 215

# --- Sheet 4: 全部类型 ---
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(2, 6).Value = 215
$ws.Cells.Item(10, 6).Value = 5091
$ws.Cells.Item(11, 6).Value = 5091
$ws.Cells.Item(12, 6).Value = 99
$ws.Cells.Item(15, 6).Value = 77
$ws.Cells.Item(16, 6).Value = 703
$ws.Cells.Item(17, 6).Value = 4929
$ws.Cells.Item(19, 6).Value = 55
$ws.Cells.Item(20, 6).Value = 73
$ws.Cells.Item(21, 6).Value = 209
$ws.Cells.Item(22, 6).Value = 97
$ws.Cells.Item(24, 6).Value = 3630
$ws.Cells.Item(25, 6).Value = 170
$ws.Cells.Item(26, 6).Value = 166
$ws.Cells.Item(27, 6).Value = 206
$ws.Cells.Item(28, 6).Value = 231
$ws.Cells.Item(29, 6).Value = 199
$ws.Cells.Item(31, 6).Value = 105
$ws.Cells.Item(36, 6).Value = 6422
$ws.Cells.Item(37, 6).Value = 1016
$ws.Cells.Item(39, 6).Value = 970
$ws.Cells.Item(40, 6).Value = 1311
$ws.Cells.Item(41, 6).Value = 153
$ws.Cells.Item(42, 6).Value = 638
$ws.Cells.Item(44, 6).Value = 2211
$ws.Cells.Item(45, 6).Value = 310
$ws.Cells.Item(48, 6).Value = 761
$ws.Cells.Item(49, 6).Value = 904
